$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------
# 1. Clear existing content
# ---------------------------------------------------------------
$ws.Cells.Clear()

# ---------------------------------------------------------------
# 2. Set new string values in the precise order needed to
#    reproduce the shared-string table ordering of the target file
# ---------------------------------------------------------------
$ws.Range("A3").Value2 = "Models"
$ws.Range("D3").Value2 = "F1-Score"
$ws.Range("A5").Value2 = "1. AltmanZ as the Only Covariate"
$ws.Range("A4").Value2 = "0. Majority Baseline"
$ws.Range("A6").Value2 = "2. All Numerical Financial Features"
$ws.Range("A9").Value2 = "5. Categorical + Numerical + NLP Features"
$ws.Range("A8").Value2 = "4. Categorical and Numerical Features"
$ws.Range("A7").Value2 = "3. NLP and Numerical Financial Features"

# Remaining header text (re-use existing shared strings)
$ws.Range("B2").Value2 = "Weighted Average"
$ws.Range("B3").Value2 = "Precision"
$ws.Range("C3").Value2 = "Recall"
$ws.Range("E3").Value2 = "Accuracy"

# ---------------------------------------------------------------
# 3. Numeric data
# ---------------------------------------------------------------
$ws.Range("E4").Value2 = 0.32

$ws.Range("B5").Value2 = 0.4
$ws.Range("C5").Value2 = 0.43
$ws.Range("D5").Value2 = 0.4
$ws.Range("E5").Value2 = 0.43

$ws.Range("B6").Value2 = 0.49
$ws.Range("C6").Value2 = 0.51
$ws.Range("D6").Value2 = 0.45
$ws.Range("E6").Value2 = 0.51

$ws.Range("B7").Value2 = 0.59
$ws.Range("C7").Value2 = 0.56999999999999995
$ws.Range("D7").Value2 = 0.56000000000000005
$ws.Range("E7").Value2 = 0.56999999999999995

$ws.Range("B8").Value2 = 0.95
$ws.Range("C8").Value2 = 0.95
$ws.Range("D8").Value2 = 0.95
$ws.Range("E8").Value2 = 0.95

$ws.Range("B9").Value2 = 0.77
$ws.Range("C9").Value2 = 0.78
$ws.Range("D9").Value2 = 0.77
$ws.Range("E9").Value2 = 0.78

# ---------------------------------------------------------------
# 4. Merge the "Weighted Average" header cell
# ---------------------------------------------------------------
$ws.Range("B2:D2").Merge()

# ---------------------------------------------------------------
# 5. Number formats (2 decimal places) for rows 4 & 5
# ---------------------------------------------------------------
$ws.Range("B4:E5").NumberFormat = "0.00"

# ---------------------------------------------------------------
# 6. Fonts - base size 12 for the whole table block, bold for the
#    "Weighted Average" merged header
# ---------------------------------------------------------------
$ws.Range("A2:E9").Font.Size = 12
$ws.Range("A2:E9").Font.Name = "Aptos Narrow"
$ws.Range("B2:D2").Font.Bold = $true

# ---------------------------------------------------------------
# 7. Center alignment for the whole block
# ---------------------------------------------------------------
$ws.Range("A2:E9").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# 8. Borders - medium weight box/grid lines
# ---------------------------------------------------------------
# Row 2 (header band) outer box + internal separators
$ws.Range("A2").Borders.Item(10).Weight = -4138   # A2 right? (not used) placeholder removed below

Write-Host "data written"
